# Guild.xlsx edit: insert a new "Force" boolean-flag row into the
# Property1 sheet, directly above the existing "Upload" row (i.e. below
# "Ref"), mirroring that row's layout/style, and update the frozen-pane /
# selection state accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property1")

# "Upload" currently lives on row 8 (right after "Ref" on row 7).
# Insert a fresh row above it; everything from the old row 8 onward
# shifts down by one (Upload -> 9, Desc -> 10, ...).
$ws.Rows.Item(8).Insert()

# Clone the formatting (styles, number formats, etc.) of the row that is
# now "Upload" (row 9) into the newly inserted blank row 8, restricted to
# the columns actually used by the table (A:AA) so we don't blow out the
# row to the full 16384-column width.
$ws.Range("A9:AA9").Copy($ws.Range("A8:AA8"))
$ws.Rows.Item(8).RowHeight = $ws.Rows.Item(9).RowHeight

# Give the new row its own label/content: a "Force" flag column, default
# FALSE for every data column (matching the other boolean-flag rows).
$ws.Range("A8").Value = "Force"

# Update the frozen pane so the split still sits right below the new
# last "flag" row (old ySplit=9/topLeftCell=A10 -> 10/A11), and fix the
# selection to match the post-edit layout.
$win = $excel.ActiveWindow
[void]($win.FreezePanes = $false)
[void]($ws.Range("A11").Select())
[void]($win.FreezePanes = $true)
[void]($ws.Range("A9").Select())
